# ooutput update 2025 august
# Refresh the generated FHIR IG "StructureDefinition-exposure-location" export:
#  - canonical base URL moved from the old GitHub shorthand repo to 2rdoc.pt
#  - regeneration timestamp bumped
#  - "Elements" sheet column widths re-settled (auto-fit) after the text change

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------------
$ws1.Range("B2").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/exposure-location"
$ws1.Range("B8").Value = "2025-08-20T10:40:04+01:00"

# --- Elements sheet --------------------------------------------------------
# Extension.url's Fixed Value mirrors the canonical URL (same underlying
# string as Metadata!B2 in the source data).
$ws2.Range("R5").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/exposure-location"
$ws2.Range("Z6").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/ValueSet/exposure-location-vs"

# Re-settle column widths (equivalent to the publisher re-running autofit
# over the "Elements" sheet once the shorter canonical URL changed the
# longest-content column). Columns whose width is capped/unaffected
# (C=10,L=12..N=14,P=16..S=19,AI=35) are left untouched.
$ws2.Columns.Item(1).ColumnWidth  = 15.666666666666666
$ws2.Columns.Item(2).ColumnWidth  = 15.666666666666666
$ws2.Columns.Item(3).ColumnWidth  = 9.0
$ws2.Columns.Item(3).Hidden       = $true
$ws2.Columns.Item(4).ColumnWidth  = 6.166666666666667
$ws2.Columns.Item(4).Hidden       = $true
$ws2.Columns.Item(5).ColumnWidth  = 4.5
$ws2.Columns.Item(6).ColumnWidth  = 3.1666666666666665
$ws2.Columns.Item(7).ColumnWidth  = 3.5
$ws2.Columns.Item(8).ColumnWidth  = 11.833333333333334
$ws2.Columns.Item(9).ColumnWidth  = 9.666666666666666
$ws2.Columns.Item(11).ColumnWidth = 13.5
$ws2.Columns.Item(15).ColumnWidth = 11.5
$ws2.Columns.Item(20).ColumnWidth = 7.0
$ws2.Columns.Item(21).ColumnWidth = 12.833333333333334
$ws2.Columns.Item(22).ColumnWidth = 13.166666666666666
$ws2.Columns.Item(23).ColumnWidth = 14.166666666666666
$ws2.Columns.Item(24).ColumnWidth = 13.833333333333334
$ws2.Columns.Item(25).ColumnWidth = 16.166666666666668
$ws2.Columns.Item(26).ColumnWidth = 57.666666666666664
$ws2.Columns.Item(27).ColumnWidth = 4.166666666666667
$ws2.Columns.Item(28).ColumnWidth = 17.166666666666668
$ws2.Columns.Item(29).ColumnWidth = 33.666666666666664
$ws2.Columns.Item(30).ColumnWidth = 12.666666666666666
$ws2.Columns.Item(31).ColumnWidth = 10.5
$ws2.Columns.Item(31).Hidden      = $true
$ws2.Columns.Item(32).ColumnWidth = 14.166666666666666
$ws2.Columns.Item(32).Hidden      = $true
$ws2.Columns.Item(33).ColumnWidth = 7.333333333333333
$ws2.Columns.Item(33).Hidden      = $true
$ws2.Columns.Item(34).ColumnWidth = 7.666666666666667
$ws2.Columns.Item(37).ColumnWidth = 18.666666666666668
